$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "52.004.93"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  +4.69%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.779.88"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  +4.83%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "342.58"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +4.65%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "115.20"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +2.10%  "

$ws.Range("E7").Value = "  +4.39%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  +4.64%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.86"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +5.24%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0857"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +4.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.01"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -0.06%  "

$ws.Range("E13").Value = "  +1.61%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.63"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +0.10%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.215.59"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  +4.89%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.773.61"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +5.23%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "51.867.81"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  +4.53%  "

$ws.Range("E18").Value = "  +1.85%  "

$ws.Range("E19").Value = "  +9.65%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.02"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +4.84%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.24"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -0.83%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0976"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +2.64%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "276.03"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +2.59%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.94"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +1.17%  "

$ws.Range("E25").Value = "  +7.30%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.66"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  +1.92%  "

$ws.Range("E27").Value = "  -0.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.16"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -0.23%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.23"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +0.82%  "

$ws.Range("E30").Value = "  +1.76%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.59"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -0.85%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.29"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +1.41%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.72"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +4.39%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0817"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -0.61%  "

$ws.Range("E35").Value = "  +0.05%  "

$ws.Range("E36").Value = "  +3.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.95"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -1.31%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.95"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +0.29%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.21"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +2.60%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0382"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +11.37%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.66"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +26.57%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.35"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  +2.82%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.116"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +3.20%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "125.95"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -2.78%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "23.12"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -2.73%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.067.10"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +0.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.31"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -0.19%  "

$ws.Range("E48").Value = "  +0.81%  "

$ws.Range("E49").Value = "  +5.42%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.83"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -1.28%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.882"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +14.62%  "
